# Deploy update: KiCad 7.0.9 -> 7.0.11+1, regenerated 2024-03-12, KiBot v1.6.4
# (pedalboard-soundcard BoM/Costs workbook)

$wb = $excel.ActiveWorkbook

$oldKicadVersion = "7.0.9-7.0.9~ubuntu23.04.1"
$newKicadVersion = "7.0.11+1"
$newCreated      = "2024-03-12 14:08:07"
$newKicost       = "KiCost® v1.1.18 + KiBot v1.6.4"

# --- "KiCad Version:" value cell on every sheet that carries it ---
$ws = $wb.Worksheets.Item("BoM")
$ws.Range("D6").Value = $newKicadVersion

$ws = $wb.Worksheets.Item("DNF")
$ws.Range("D6").Value = $newKicadVersion

$ws = $wb.Worksheets.Item("Costs")
$ws.Range("E6").Value = $newKicadVersion

$ws = $wb.Worksheets.Item("Costs (DNF)")
$ws.Range("E6").Value = $newKicadVersion

# --- "Created:" timestamp + KiCost/KiBot version caption ---
$ws = $wb.Worksheets.Item("Costs")
$ws.Range("B30").Value = $newCreated
$ws.Range("A31").Value = $newKicost

$ws = $wb.Worksheets.Item("Costs (DNF)")
$ws.Range("B13").Value = $newCreated
$ws.Range("A14").Value = $newKicost

# --- Column width tweaks (narrower "Datasheet" / unit-cost columns) ---
# Excel quantises ColumnWidth to sixths of a character (pixel grid at MDW=6
# in this runtime), so we pick the ColumnWidth whose resulting stored width
# lands closest to the authored target width.
$ws = $wb.Worksheets.Item("DNF")
$ws.Columns.Item(4).ColumnWidth = 19.833333333333332   # -> width ~20.71 (was 25.71)

$ws = $wb.Worksheets.Item("Costs")
$ws.Columns.Item(5).ColumnWidth = 20.833333333333332   # -> width ~21.71 (was 26.71)

$ws = $wb.Worksheets.Item("Costs (DNF)")
$ws.Columns.Item(5).ColumnWidth = 20.833333333333332   # -> width ~21.71 (was 26.71)
